# New crime data collected - weekly CompStat update for 116th Precinct
# Report period moves forward one week (6/16-6/22 -> 6/23-6/29), volume number 25 -> 26,
# and the weekly/28-day/YTD/2-year crime statistics table (rows 15-30) is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume/Number and report date range ---
$ws.Range("A8").Value = "Volume 32   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# --- Row 15 (Murder) ---
# C15 goes from a number (2) to the textual "0" placeholder, matching D15's style.
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 16 (Rape) ---
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -90.909090909090
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = -38.636363636363

# --- Row 17 (Robbery) ---
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -10.526315789473
$ws.Range("I17").Value = 124
$ws.Range("J17").Value = 115
$ws.Range("K17").Value = 7.826086956521
$ws.Range("L17").Value = -2.362204724409

# --- Row 18 (Fel. Assault) ---
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = -42.5
$ws.Range("L18").Value = -34.285714285714

# --- Row 19 (Burglary) ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 109
$ws.Range("J19").Value = 131
$ws.Range("K19").Value = -16.793893129771
$ws.Range("L19").Value = -18.656716417910

# --- Row 20 (Gr. Larceny) ---
$ws.Range("C20").Value = 7
# D20 and E20 go from text placeholders ("0" / "***.*") to real numbers; adopt the
# same numeric styles already used by the other numeric cells in this table (s=14 / s=15).
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 6
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = 16.666666666666
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 15.384615384615
$ws.Range("I20").Value = 94
$ws.Range("J20").Value = 79
$ws.Range("K20").Value = 18.987341772151
$ws.Range("L20").Value = 1.075268817204

# --- Row 21 (G.L.A.) ---
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -22.727272727272
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = -17.910447761194
$ws.Range("I21").Value = 395
$ws.Range("J21").Value = 415
$ws.Range("K21").Value = -4.819277108433
$ws.Range("L21").Value = -9.403669724770

# --- Row 24 (Transit) ---
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 10
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = -19.512195121951
$ws.Range("I24").Value = 236
$ws.Range("J24").Value = 292
$ws.Range("K24").Value = -19.178082191780
$ws.Range("L24").Value = -9.578544061302

# --- Row 25 (Housing) ---
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -68.75
$ws.Range("J25").Value = 71
$ws.Range("K25").Value = -22.535211267605
$ws.Range("L25").Value = 1.851851851851

# --- Row 26 (Petit Larceny) ---
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 18.181818181818
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -2.380952380952
$ws.Range("I26").Value = 217
$ws.Range("J26").Value = 215
$ws.Range("K26").Value = 0.930232558139
$ws.Range("L26").Value = 33.128834355828

# --- Row 27 (Retail Theft) ---
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = 91.666666666666

# --- Row 28 (Misd. Assault) ---
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0

# --- Row 29 (UCR Rape*) ---
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 30 (Other Sex Crimes) ---
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
